# Apply "Latest TestNG changes are included" edits to the workbook.
# Target: xl/worksheets/sheet1.xml (DemoWebShop) is updated with new
# login credentials, the UserName column becomes Email, hyperlinks are
# removed, and the active sheet/selection moves from OrangeHRM (C29 /
# O3) to DemoWebShop (C4).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DemoWebShop")
$ws2 = $wb.Worksheets.Item("OrangeHRM")

# Remove every hyperlink that lives on the DemoWebShop sheet (C2, D2,
# C3, D3) - deleting the whole collection drops the <hyperlinks> node
# and the now-unused relationship entries.
$ws1.Hyperlinks.Delete()

# Row 2: new login e-mail / password pair
$ws1.Range("C2").Value = "log-in123@gmail.com"
$ws1.Range("D2").Value = "login@123"

# Row header: "UserName" -> "Email"
$ws1.Range("C1").Value = "Email"

# Row 3: credentials removed entirely (now blank, but formatted like C2/D2)
$ws1.Range("C3").ClearContents()
$ws1.Range("D3").ClearContents()

# A2 / A3 become text ("1"/"2" with a quote-prefix) instead of numbers 1/2
$ws1.Range("A2").Value = "'1"
$ws1.Range("A3").Value = "'2"

# Make DemoWebShop the active sheet/tab with C4 selected (was OrangeHRM
# with O3 selected before).
$ws1.Activate()
$ws1.Range("C4").Select()
